$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.062008445517327
$ws.Range("C2").Value = 0.062008445517327
$ws.Range("D2").Value = 3.27833851059219
$ws.Range("F2").Value = 0.0709

$ws.Range("B3").Value = 2.76153088409848
$ws.Range("C3").Value = 0.018914595096565
